$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.317.25"
$ws.Range("E2").Value = "  -4.91%  "
$ws.Range("D3").Value = "3.256.55"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.22"
$ws.Range("E5").Value = "  -5.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.77"
$ws.Range("E6").Value = "  -12.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.248.95"
$ws.Range("E8").Value = "  -7.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -10.92%  "
$ws.Range("E10").Value = "  -13.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.81"
$ws.Range("E11").Value = "  -4.33%  "
$ws.Range("E12").Value = "  -13.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.66"
$ws.Range("E13").Value = "  -16.90%  "
$ws.Range("E14").Value = "  -11.61%  "
$ws.Range("D15").Value = "3.773.15"
$ws.Range("E15").Value = "  -7.80%  "
$ws.Range("D16").Value = "67.335.51"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "549.94"
$ws.Range("E17").Value = "  -9.48%  "
$ws.Range("D18").Value = "3.255.35"
$ws.Range("E18").Value = "  -7.70%  "
$ws.Range("E19").Value = "  -13.50%  "
$ws.Range("E20").Value = "  -5.91%  "
$ws.Range("E21").Value = "  -14.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.770"
$ws.Range("E22").Value = "  -13.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.97"
$ws.Range("E23").Value = "  -12.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.82"
$ws.Range("E24").Value = "  -12.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.61"
$ws.Range("E25").Value = "  -13.17%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.23"
$ws.Range("E27").Value = "  -14.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  -10.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.61"
$ws.Range("E29").Value = "  -12.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("E30").Value = "  -16.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  -12.45%  "
$ws.Range("E32").Value = "  -11.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "548.66"
$ws.Range("E33").Value = "  -13.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.68"
$ws.Range("E34").Value = "  -17.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.79"
$ws.Range("E35").Value = "  -14.94%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0451"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.67"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0860"
$ws.Range("E39").Value = "  -13.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.26"
$ws.Range("E40").Value = "  -14.55%  "
$ws.Range("E41").Value = "  -11.93%  "
$ws.Range("D42").Value = "2.937.55"
$ws.Range("E42").Value = "  -12.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  -23.51%  "
$ws.Range("E44").Value = "  -15.61%  "
$ws.Range("D45").Value = "0.0₃0587"
$ws.Range("E45").Value = "  -19.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.64"
$ws.Range("E46").Value = "  -17.01%  "
$ws.Range("E47").Value = "  -20.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.16"
$ws.Range("E48").Value = "  -15.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.19"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.115"
$ws.Range("E51").Value = "  -12.21%  "
